$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = $ws.UsedRange.Rows.Count }

for ($r = 2; $r -le $lastRow; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)

    $bVal = $bCell.Value()
    $cVal = $cCell.Value()

    if ([string]::IsNullOrEmpty($bVal) -or [string]::IsNullOrEmpty($cVal)) {
        continue
    }

    # e.g. "club-sports" -> "club" ; "uil-sports" -> "uil"
    $bPrefix = $bVal.Split('-')[0]

    # e.g. "Basketball-Boys" -> "boys"
    $cParts = $cVal.Split('-')
    $cSuffix = $cParts[1].ToLower()

    $newB = "sports_" + $bPrefix + "_" + $cSuffix
    $newC = $cParts[0]

    $bCell.Value = $newB
    $cCell.Value = $newC
}
